$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$newValues = @(
    "48×23=1104",
    "90×30=2700",
    "97×76=7372",
    "50×39=1950",
    "19×29=551",
    "100×95=9500",
    "52×41=2132",
    "46×87=4002",
    "82×84=6888",
    "64×65=4160",
    "95×30=2850",
    "99×72=7128",
    "58×24=1392",
    "97×32=3104",
    "80×68=5440",
    "71×50=3550",
    "76×40=3040",
    "45×30=1350",
    "39×68=2652",
    "35×75=2625",
    "52×80=4160",
    "79×17=1343",
    "84×21=1764",
    "11×43=473",
    "72×54=3888",
    "99×18=1782",
    "17×84=1428",
    "84×56=4704",
    "12×16=192",
    "75×21=1575",
    "47×87=4089",
    "67×19=1273",
    "48×62=2976",
    "73×28=2044",
    "98×86=8428",
    "45×39=1755",
    "79×72=5688",
    "31×39=1209",
    "68×25=1700",
    "23×98=2254",
    "62×72=4464",
    "45×74=3330",
    "46×12=552",
    "38×67=2546",
    "70×16=1120",
    "54×60=3240",
    "29×27=783",
    "50×20=1000",
    "47×14=658",
    "67×48=3216",
    "13×66=858",
    "79×87=6873",
    "80×22=1760",
    "70×44=3080",
    "50×28=1400",
    "68×16=1088",
    "93×14=1302",
    "16×93=1488",
    "78×20=1560",
    "87×72=6264",
    "58×48=2784",
    "100×98=9800",
    "40×97=3880",
    "20×27=540",
    "29×20=580",
    "36×82=2952",
    "75×50=3750",
    "64×87=5568",
    "39×30=1170",
    "61×74=4514",
    "66×69=4554",
    "42×10=420",
    "69×56=3864",
    "12×28=336",
    "41×49=2009",
    "80×36=2880",
    "65×76=4940",
    "25×28=700",
    "90×69=6210",
    "46×53=2438",
    "62×12=744",
    "72×71=5112",
    "98×39=3822",
    "33×53=1749",
    "48×97=4656",
    "19×90=1710",
    "61×98=5978",
    "32×77=2464",
    "69×77=5313",
    "51×17=867",
    "82×83=6806",
    "83×35=2905",
    "86×58=4988",
    "78×78=6084",
    "20×43=860",
    "14×10=140",
    "58×78=4524",
    "11×83=913",
    "58×69=4002",
    "48×84=4032"
)

$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}
Write-Host "Updated" $idx "cells"
